# Initial TTS and TTLE calibration
# Update the Transportation Technology Logit Exponents on the TTLE sheet
# from -3 to -2 for both Passenger and Freight columns (B2:C7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TTLE")

$ws.Range("B2:C7").Value = -2

# Make the TTLE sheet the active sheet with B2:C7 selected, matching the
# saved view state captured in the workbook after calibration.
$ws.Activate()
$ws.Range("B2:C7").Select()
